$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.79549166666667
$ws.Range("H2").Value = 38.386475
$ws.Range("I2").Value = 0.5145949251267348
$ws.Range("J2").Value = 0.5145949251267348
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.126378
$ws.Range("N2").Value = 12.379134
$ws.Range("O2").Value = 0.1457327627461222
$ws.Range("P2").Value = 0.1457327627461223
$ws.Range("Q2").Value = 52.79903531251666
$ws.Range("R2").Value = 475.19131781265
$ws.Range("S2").Value = 0.07499334013385298
$ws.Range("T2").Value = 0.07499334013385299

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.79549166666667
$ws.Range("H3").Value = 38.386475
$ws.Range("I3").Value = 0.5145949251267348
$ws.Range("J3").Value = 0.5145949251267348
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.39252366666667
$ws.Range("N3").Value = 34.177571
$ws.Range("O3").Value = 0.4023538194014015
$ws.Range("P3").Value = 0.4023538194014014
$ws.Range("Q3").Value = 145.7729416391361
$ws.Range("R3").Value = 1311.956474752225
$ws.Range("S3").Value = 0.20704923356932
$ws.Range("T3").Value = 0.2070492335693199

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.79549166666667
$ws.Range("H4").Value = 38.386475
$ws.Range("I4").Value = 0.5145949251267348
$ws.Range("J4").Value = 0.5145949251267348
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 12.79578833333333
$ws.Range("N4").Value = 38.387365
$ws.Range("O4").Value = 0.4519134178524764
$ws.Range("P4").Value = 0.4519134178524763
$ws.Range("Q4").Value = 163.7284029875972
$ws.Range("R4").Value = 1473.555626888375
$ws.Range("S4").Value = 0.2325523514235619
$ws.Range("T4").Value = 0.2325523514235618

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.803896
$ws.Range("H5").Value = 29.411688
$ws.Range("I5").Value = 0.3942822409249843
$ws.Range("J5").Value = 0.3942822409249843
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.126378
$ws.Range("N5").Value = 12.379134
$ws.Range("O5").Value = 0.1457327627461222
$ws.Range("P5").Value = 0.1457327627461223
$ws.Range("Q5").Value = 40.454580768688
$ws.Range("R5").Value = 364.091226918192
$ws.Range("S5").Value = 0.05745984027173014
$ws.Range("T5").Value = 0.05745984027173015

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 9.803896
$ws.Range("H6").Value = 29.411688
$ws.Range("I6").Value = 0.3942822409249843
$ws.Range("J6").Value = 0.3942822409249843
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 11.39252366666667
$ws.Range("N6").Value = 34.177571
$ws.Range("O6").Value = 0.4023538194014015
$ws.Range("P6").Value = 0.4023538194014014
$ws.Range("Q6").Value = 111.6911172055387
$ws.Range("R6").Value = 1005.220054849848
$ws.Range("S6").Value = 0.158640965558311
$ws.Range("T6").Value = 0.1586409655583109

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 9.803896
$ws.Range("H7").Value = 29.411688
$ws.Range("I7").Value = 0.3942822409249843
$ws.Range("J7").Value = 0.3942822409249843
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 12.79578833333333
$ws.Range("N7").Value = 38.387365
$ws.Range("O7").Value = 0.4519134178524764
$ws.Range("P7").Value = 0.4519134178524763
$ws.Range("Q7").Value = 125.4485780580133
$ws.Range("R7").Value = 1129.03720252212
$ws.Range("S7").Value = 0.1781814350949432
$ws.Range("T7").Value = 0.1781814350949431

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.265785
$ws.Range("H8").Value = 6.797355
$ws.Range("I8").Value = 0.09112283394828093
$ws.Range("J8").Value = 0.09112283394828093
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.126378
$ws.Range("N8").Value = 12.379134
$ws.Range("O8").Value = 0.1457327627461222
$ws.Range("P8").Value = 0.1457327627461223
$ws.Range("Q8").Value = 9.34948537673
$ws.Range("R8").Value = 84.14536839057001
$ws.Range("S8").Value = 0.01327958234053912
$ws.Range("T8").Value = 0.01327958234053912

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.265785
$ws.Range("H9").Value = 6.797355
$ws.Range("I9").Value = 0.09112283394828093
$ws.Range("J9").Value = 0.09112283394828093
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 11.39252366666667
$ws.Range("N9").Value = 34.177571
$ws.Range("O9").Value = 0.4023538194014015
$ws.Range("P9").Value = 0.4023538194014014
$ws.Range("Q9").Value = 25.81300923607834
$ws.Range("R9").Value = 232.317083124705
$ws.Range("S9").Value = 0.03666362027377052
$ws.Range("T9").Value = 0.03666362027377051

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.265785
$ws.Range("H10").Value = 6.797355
$ws.Range("I10").Value = 0.09112283394828093
$ws.Range("J10").Value = 0.09112283394828093
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 12.79578833333333
$ws.Range("N10").Value = 38.387365
$ws.Range("O10").Value = 0.4519134178524764
$ws.Range("P10").Value = 0.4519134178524763
$ws.Range("Q10").Value = 28.99250526884166
$ws.Range("R10").Value = 260.932547419575
$ws.Range("S10").Value = 0.0411796313339713
$ws.Range("T10").Value = 0.04117963133397129

